# Commit: "Doing Updates for Financials"
# A new fiscal-period column is inserted before column D (the existing D:K
# columns, and their data, shift one column right to E:L); the new column D
# is then populated with the figures for the newest reporting period.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HAL")

$ws.Columns("D:D").Insert()

# Exact (escaped) format codes from the workbook, hard-coded so re-applying them
# reuses the existing style indexes instead of minting near-duplicate numFmts
# (the NumberFormat getter normalises away the backslash escapes on read-back).
$dateFormat = '[$-409]d\-mmm\-yy;@'
$numberFormat = '#,##0'

# Period-ending date header cells (row 7/38/80) -> 2019-01-01 (serial 43465)
$dateValues = @{
    7 = 43465
    38 = 43465
    80 = 43465
}

# Plain numeric data cells for the new period
$numValues = @{
    8 = 23995000
    9 = 21009000
    10 = 2986000
    13 = 0
    14 = 265000
    15 = 0
    17 = 21528000
    18 = 2467000
    20 = -55000
    21 = 4018000
    22 = 598000
    23 = 1814000
    24 = 204000
    25 = 0
    26 = 1610000
    27 = 1609000
    28 = 0
    29 = 47000
    30 = 0
    31 = 0
    32 = 55000
    33 = 1656000
    34 = 0
    35 = 1656000
    41 = 2008000
    43 = 5234000
    44 = 3028000
    45 = 881000
    46 = 11151000
    48 = 8961000
    49 = 2825000
    50 = 0
    51 = 0
    52 = 3045000
    53 = 0
    54 = 25982000
    57 = 3018000
    58 = 36000
    59 = 1748000
    60 = 4802000
    61 = 10421000
    62 = 1215000
    63 = 0
    64 = 0
    65 = 0
    66 = 16460000
    68 = 0
    69 = 0
    70 = 0
    71 = 0
    72 = 13739000
    73 = 0
    74 = 0
    75 = 0
    76 = 9522000
    77 = 0
    81 = 1656000
    83 = 1606000
    84 = 0
    85 = 0
    86 = 0
    87 = 0
    88 = 0
    89 = 3157000
    91 = -2026000
    92 = 0
    93 = 0
    94 = -1993000
    96 = -630000
    97 = 0
    98 = 0
    99 = 0
    100 = -1419000
    101 = -74000
    102 = -329000
}

# Rows whose new cell is the literal "NA" placeholder (shared string), matching
# the neighbouring cells in the same row
$naRows = @(12, 42, 47)

# Rows that stay blank in every period column
$emptyRows = @(11, 16, 19, 39, 40, 55, 56, 67, 82, 90, 95)

foreach ($r in $dateValues.Keys) {
    $cell = $ws.Cells.Item($r, 4)
    $cell.NumberFormat = $dateFormat
    $cell.Value = $dateValues[$r]
}

foreach ($r in $numValues.Keys) {
    $cell = $ws.Cells.Item($r, 4)
    $cell.NumberFormat = $numberFormat
    $cell.Value = $numValues[$r]
}

foreach ($r in $naRows) {
    $cell = $ws.Cells.Item($r, 4)
    $cell.NumberFormat = $numberFormat
    $cell.Value = "NA"
}

foreach ($r in $emptyRows) {
    $cell = $ws.Cells.Item($r, 4)
    $cell.NumberFormat = $numberFormat
}
